# TC05 - Search Product
# Adds a new "TC05" worksheet (mirroring the baseUrl/email header pattern
# used by the other test-case sheets) and updates the sheet selections that
# shifted as part of this change.

$wb = $excel.ActiveWorkbook

# --- TC01: re-apply focus on K2 -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$null = $ws1.Select()
$ws1Sel = $ws1.Range("K2")
$null = $ws1Sel.Select()

# --- TC02: selection grows from A1 to A1:A2 -------------------------------
$ws2 = $wb.Worksheets.Item(2)
$null = $ws2.Select()
$ws2Sel = $ws2.Range("A1:A2")
$null = $ws2Sel.Select()

# --- TC03: selection grows from B2 to A1:A2 -------------------------------
$ws3 = $wb.Worksheets.Item(3)
$null = $ws3.Select()
$ws3Sel = $ws3.Range("A1:A2")
$null = $ws3Sel.Select()

# --- TC04: selection moves from B3 to A18 (sheet no longer the active tab)
$ws4 = $wb.Worksheets.Item(4)
$null = $ws4.Select()
$ws4Sel = $ws4.Range("A18")
$null = $ws4Sel.Select()

# --- Add the new TC05 sheet after the last existing sheet -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "TC05"

$ws5.Range("A1").Value = "baseUrl"
$ws5.Range("A2").Value = "http://automationpractice.com/index.php"

# Match the column A width used by the sibling TC01-TC04 sheets.
$ws5Col = $ws5.Columns.Item(1)
$ws5Col.ColumnWidth = 33

# Match the print/page setup used by the sibling TC01-TC04 sheets.
$ps5 = $ws5.PageSetup
$ps5.LeftMargin = 56.7
$ps5.RightMargin = 56.7
$ps5.TopMargin = 75.8
$ps5.BottomMargin = 75.8
$ps5.HeaderMargin = 56.7
$ps5.FooterMargin = 56.7
$ps5.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ps5.CenterFooter = '&"Times New Roman,Normal"&12Página &P'
$ps5.PrintHeadings = $false
$ps5.PrintGridlines = $false
$ps5.CenterHorizontally = $false
$ps5.CenterVertically = $false
$ps5.PaperSize = 1
$ps5.Zoom = 100
$ps5.FitToPagesWide = 1
$ps5.FitToPagesTall = 1
$ps5.Orientation = 1
$ps5.BlackAndWhite = $false
$ps5.Draft = $false

# TC05 becomes the active sheet/selection (A1:A2), matching the new tab
# being selected last.
$null = $ws5.Select()
$ws5Sel = $ws5.Range("A1:A2")
$null = $ws5Sel.Select()
